$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.240.35'
$ws.Range("E2").Value = '  -0.10%  '

$ws.Range("D3").Value = '3.343.89'
$ws.Range("E3").Value = '  +0.38%  '

$ws.Range("D4").Value = "'0.997"
$ws.Range("E4").Value = '  -0.42%  '

$ws.Range("D5").Value = "'584.61"
$ws.Range("E5").Value = '  +4.15%  '

$ws.Range("D6").Value = "'185.52"
$ws.Range("E6").Value = '  -1.79%  '

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").Value = '3.340.83'
$ws.Range("E8").Value = '  +0.50%  '

$ws.Range("E9").Value = '  -2.47%  '

$ws.Range("D10").Value = "'0.181"
$ws.Range("E10").Value = '  -1.83%  '

$ws.Range("D11").Value = "'0.583"
$ws.Range("E11").Value = '  -1.39%  '

$ws.Range("D12").Value = "'47.02"

$ws.Range("E13").Value = '  -1.41%  '

$ws.Range("D14").Value = "'668.66"
$ws.Range("E14").Value = '  +10.21%  '

$ws.Range("D15").Value = '3.879.42'
$ws.Range("E15").Value = '  +0.40%  '

$ws.Range("D16").Value = "'8.52"
$ws.Range("E16").Value = '  -2.11%  '

$ws.Range("D17").Value = '66.450.79'
$ws.Range("E17").Value = '  +0.10%  '

$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").Value = "'0.118"
$ws.Range("E18").Value = '  -0.67%  '

$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").Value = "'17.90"
$ws.Range("E19").Value = '  -1.14%  '

$ws.Range("D20").Value = '3.343.02'
$ws.Range("E20").Value = '  +0.09%  '

$ws.Range("D21").Value = "'11.08"
$ws.Range("E21").Value = '  -0.80%  '

$ws.Range("D22").Value = "'0.898"
$ws.Range("E22").Value = '  -2.02%  '

$ws.Range("D23").Value = "'17.75"
$ws.Range("E23").Value = '  -4.55%  '

$ws.Range("D24").Value = "'101.45"
$ws.Range("E24").Value = '  +0.79%  '

$ws.Range("E25").Value = '  -1.95%  '

$ws.Range("D26").Value = "'3.98"

$ws.Range("D27").Value = "'2.78"
$ws.Range("E27").Value = '  +0.08%  '

$ws.Range("D28").Value = "'9.44"
$ws.Range("E28").Value = '  -2.90%  '

$ws.Range("D29").Value = "'32.16"
$ws.Range("E29").Value = '  +5.20%  '

$ws.Range("D30").Value = "'8.51"

$ws.Range("D31").Value = "'6.83"
$ws.Range("E31").Value = '  +0.18%  '

$ws.Range("D32").Value = "'610.90"
$ws.Range("E32").Value = '  +4.90%  '

$ws.Range("D33").Value = "'3.89"
$ws.Range("E33").Value = '  +0.26%  '

$ws.Range("D34").Value = "'11.11"
$ws.Range("E34").Value = '  -0.64%  '

$ws.Range("D35").Value = '3.859.43'
$ws.Range("E35").Value = '  +3.61%  '

$ws.Range("D36").Value = "'0.105"
$ws.Range("E36").Value = '  -0.73%  '

$ws.Range("E37").Value = '  +0.01%  '

$ws.Range("D38").Value = "'56.24"

$ws.Range("D39").Value = "'0.128"
$ws.Range("E39").Value = '  -2.85%  '

$ws.Range("B40").Value = 'PEPE'
$ws.Range("C40").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D40").Value = '0.0₃0702'
$ws.Range("E40").Value = '  -4.17%  '

$ws.Range("B41").Value = 'Fetch.AI'
$ws.Range("C41").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D41").Value = "'2.67"
$ws.Range("E41").Value = '  -1.26%  '

$ws.Range("D42").Value = "'32.93"
$ws.Range("E42").Value = '  -3.37%  '

$ws.Range("D43").Value = "'3.19"
$ws.Range("E43").Value = '  -3.42%  '

$ws.Range("E44").Value = '  +1.60%  '

$ws.Range("D45").Value = "'0.337"
$ws.Range("E45").Value = '  -2.53%  '

$ws.Range("D46").Value = "'0.0417"
$ws.Range("E46").Value = '  -2.19%  '

$ws.Range("D47").Value = "'2.99"
$ws.Range("E47").Value = '  -13.83%  '

$ws.Range("B49").Value = 'FirstDigitalUSD'
$ws.Range("C49").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D49").Value = "'1.00"
$ws.Range("E49").Value = '  +0.41%  '

$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D50").Value = "'2.55"
$ws.Range("E50").Value = '  -2.37%  '

$ws.Range("E51").Value = '  +1.90%  '
